$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are written as text, preserving exact formatting
# (e.g. trailing zeros, thousand-dot separators) instead of being auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.468.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.932.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9977"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9988"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4757"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2876"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06606"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.27"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "108.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +27.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.921.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.62%  "

$ws.Range("E13").Value = "  +2.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.174"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6636"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "311.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +26.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.472.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9996"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007559"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.177.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.331"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9990"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.317"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.88%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.298"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.82%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.049"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.34%  "

$ws.Range("E29").Value = "  +8.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.361"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.121"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.942"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05031"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7445"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.156"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.757"
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01974"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.052"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8797"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.802"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9987"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4188"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.319"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.323"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1209"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05625"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3871"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.06%  "
